$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.406.34"
$ws.Range("D3").Value = "1.689.56"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").Value = "'219.16"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").Value = "'0.5518"
$ws.Range("E6").Value = "  +8.39%  "
$ws.Range("D7").Value = "'1.010"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("D8").Value = "'0.2711"
$ws.Range("E8").Value = "  +2.15%  "
$ws.Range("D9").Value = "'0.06489"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").Value = "'22.11"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").Value = "'0.07602"
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("D12").Value = "'4.568"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").Value = "1.690.95"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").Value = "'0.5824"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "'0.000008470"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "'65.42"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").Value = "26.538.58"
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("D18").Value = "'4.947"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "'1.010"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("D21").Value = "'190.61"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'6.258"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").Value = "'1.011"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").Value = "'149.57"
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("D25").Value = "'0.1317"
$ws.Range("E25").Value = "  +9.75%  "
$ws.Range("D26").Value = "'7.943"
$ws.Range("E26").Value = "  +4.34%  "
$ws.Range("D27").Value = "'15.83"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").Value = "'0.06319"
$ws.Range("E28").Value = "  -4.50%  "
$ws.Range("D29").Value = "'1.410"
$ws.Range("E29").Value = "  +5.60%  "
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").Value = "'3.592"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("D33").Value = "'1.676"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").Value = "'1.047"
$ws.Range("E34").Value = "  +2.98%  "
$ws.Range("D35").Value = "'0.6251"
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("D36").Value = "'2.410"
$ws.Range("D37").Value = "'2.722"
$ws.Range("E37").Value = "  +1.46%  "
$ws.Range("D38").Value = "'6.247"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").Value = "'0.01645"
$ws.Range("E39").Value = "  +3.34%  "
$ws.Range("D40").Value = "1.120.79"
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("D41").Value = "'0.8834"
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("D42").Value = "'1.016"
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").Value = "'100.90"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "1.841.93"
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("D45").Value = "'0.00000000112"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "'57.57"
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("D47").Value = "'8.214"
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "'0.05281"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").Value = "'0.4304"
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("E51").Value = "  +1.42%  "
